# Updated cryptos list - refresh prices / volume percentages.
# Row 35/36 also swap coin identity (HuobiToken moves to row 35, Maker to row 36)
# along with their own refreshed price/volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.848.37'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '1.643.02'
$ws.Range("E3").Value = '  +0.04%  '

$ws.Range("E4").Value = '  -0.74%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '216.54'
$ws.Range("E5").Value = '  -0.85%  '

$ws.Range("E6").Value = '  +0.98%  '

$ws.Range("E7").Value = '  -0.73%  '

$ws.Range("E8").Value = '  +1.09%  '

$ws.Range("E9").Value = '  -0.43%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.84'
$ws.Range("E10").Value = '  +4.00%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0844'
$ws.Range("E11").Value = '  -0.23%  '

$ws.Range("D12").Value = '1.873.23'
$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("D13").Value = '1.642.47'
$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("E15").Value = '  +0.63%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '66.35'
$ws.Range("E16").Value = '  +2.56%  '

$ws.Range("D17").Value = '26.871.35'
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("E18").Value = '  +0.80%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '217.61'
$ws.Range("E19").Value = '  +2.99%  '

$ws.Range("E20").Value = '  -0.80%  '

$ws.Range("E22").Value = '  +6.83%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.43'
$ws.Range("E23").Value = '  +5.20%  '

$ws.Range("E24").Value = '  -0.85%  '

$ws.Range("E25").Value = '  -1.45%  '

$ws.Range("E26").Value = '  -0.82%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.36'
$ws.Range("E27").Value = '  +4.10%  '

$ws.Range("E28").Value = '  +0.46%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.84'
$ws.Range("E29").Value = '  +1.68%  '

$ws.Range("E30").Value = '  +1.99%  '

$ws.Range("E31").Value = '  -0.36%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.36'
$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("E33").Value = '  +0.48%  '

$ws.Range("E34").Value = '  +2.08%  '

# Rows 35/36: HuobiToken and Maker swap positions, each with new figures.
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.44'
$ws.Range("E35").Value = '  -0.67%  '

$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '1.244.26'
$ws.Range("E36").Value = '  -2.52%  '

$ws.Range("E37").Value = '  -0.26%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.537'
$ws.Range("E38").Value = '  +1.96%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.834'
$ws.Range("E39").Value = '  +3.33%  '

$ws.Range("E40").Value = '  -0.75%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.808'
$ws.Range("E41").Value = '  +0.43%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.36'
$ws.Range("E42").Value = '  +1.73%  '

$ws.Range("D43").Value = '1.785.73'
$ws.Range("E43").Value = '  +0.21%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.10'
$ws.Range("E44").Value = '  -4.02%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '60.96'
$ws.Range("E45").Value = '  +1.39%  '

$ws.Range("E46").Value = '  -0.02%  '

$ws.Range("E47").Value = '  +0.43%  '

$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("E49").Value = '  -1.21%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0971'
$ws.Range("E50").Value = '  +1.20%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.54'
$ws.Range("E51").Value = '  +0.35%  '
